# Fixed #355 Added rotate(MImage) service.
# Update the m:asImage() field's instrText runs:
#   'dh1.gif'.asImage().resize(0.2)  ->  'Mona_Lisa.jpg'.asImage().resize(0.5)
# The field instruction text lives in w:instrText runs split across many
# <w:r> elements; Word's Field.Code setter only prepends a new run instead
# of rewriting the existing ones, so we rebuild the whole host paragraph's
# OOXML (fldChar begin/end, bookmark, and all instrText runs) via
# Range.InsertXML, which is the supported way to replace a Range's raw
# contents with exact WordprocessingML.

$d = $word.ActiveDocument

# Locate the paragraph that hosts the field (robust to its index).
$field = $d.Fields.Item(1)
$fieldStart = $field.Code.Start
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($fieldStart -ge $p.Range.Start -and $fieldStart -lt $p.Range.End) {
        $targetPara = $p
        break
    }
}

$newParaXml = '<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="00DE6D5A"><w:instrText>m</w:instrText></w:r><w:r><w:instrText>:</w:instrText></w:r><w:r w:rsidR="004B598D"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:instrText>''</w:instrText></w:r><w:r w:rsidR="00321AA1"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:instrText>Mona_Lisa</w:instrText></w:r><w:r w:rsidR="004B598D"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:instrText>.</w:instrText></w:r><w:r w:rsidR="004B598D"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:instrText>jpg</w:instrText></w:r><w:r w:rsidR="004B598D"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:instrText>''.asImage()</w:instrText></w:r><w:r w:rsidR="00484D7C"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:instrText>.resize(0.</w:instrText></w:r><w:r w:rsidR="001E6881"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:instrText>5</w:instrText></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="00FB48D7"><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:instrText>)</w:instrText></w:r><w:r><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p>'

$flatOpc = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $newParaXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetPara.Range.InsertXML($flatOpc)
